# Add a new "DP on Matrix" entry (Leetcode - 1277) and a new
# "Stock Graph Pattern" entry (Leetcode - 55) to the pattern tracker sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New DP problem under the "DP on Matrix" column (E) - first new row (row 4).
$ws.Range("E4").Value = "Leetcode - 1277"

# New problem under the "Stock Graph Pattern" column (I), continuing row 3.
$ws.Range("I3").Value = "Leetcode - 55"

# Leave the cursor where the author last clicked after typing the entries.
$null = $ws.Range("C5").Select()
